$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44671

# Row 3
$ws.Cells.Item(3, 4).Value = 44965
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 100
$ws.Cells.Item(3, 14).Value = 34000
$ws.Cells.Item(3, 15).Value = 35000
$ws.Cells.Item(3, 16).Value = 34600
$ws.Cells.Item(3, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(3, 19).Value = 1922
$ws.Cells.Item(3, 20).Value = 18

# Row 4
$ws.Cells.Item(4, 4).Value = 44965
$ws.Cells.Item(4, 13).Value = 120
$ws.Cells.Item(4, 14).Value = 32000
$ws.Cells.Item(4, 15).Value = 33000
$ws.Cells.Item(4, 16).Value = 32333
$ws.Cells.Item(4, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(4, 19).Value = 1796
$ws.Cells.Item(4, 20).Value = 18

# Row 5
$ws.Cells.Item(5, 4).Value = 44643
$ws.Cells.Item(5, 14).Value = 28000
$ws.Cells.Item(5, 15).Value = 30000
$ws.Cells.Item(5, 16).Value = 29000
$ws.Cells.Item(5, 19).Value = 1450

# Row 6
$ws.Cells.Item(6, 4).Value = 44679
$ws.Cells.Item(6, 13).Value = 200

# Row 7
$ws.Cells.Item(7, 4).Value = 44679
$ws.Cells.Item(7, 12).Value = "Tercera"
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 24000
$ws.Cells.Item(7, 15).Value = 25000
$ws.Cells.Item(7, 16).Value = 24500
$ws.Cells.Item(7, 19).Value = 1225

# Row 8
$ws.Cells.Item(8, 4).Value = 45021
$ws.Cells.Item(8, 13).Value = 250
$ws.Cells.Item(8, 14).Value = 22000
$ws.Cells.Item(8, 15).Value = 23000
$ws.Cells.Item(8, 16).Value = 22500
$ws.Cells.Item(8, 19).Value = 1125

# Row 9
$ws.Cells.Item(9, 4).Value = 44993
$ws.Cells.Item(9, 13).Value = 130
$ws.Cells.Item(9, 14).Value = 25000
$ws.Cells.Item(9, 15).Value = 26000
$ws.Cells.Item(9, 16).Value = 25462
$ws.Cells.Item(9, 19).Value = 1273

# Row 10
$ws.Cells.Item(10, 4).Value = 44650
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 160
$ws.Cells.Item(10, 14).Value = 31000
$ws.Cells.Item(10, 15).Value = 32000
$ws.Cells.Item(10, 16).Value = 31500
$ws.Cells.Item(10, 19).Value = 1575

# Row 11
$ws.Cells.Item(11, 4).Value = 44650
$ws.Cells.Item(11, 13).Value = 250
$ws.Cells.Item(11, 14).Value = 29000
$ws.Cells.Item(11, 15).Value = 30000
$ws.Cells.Item(11, 16).Value = 29500
$ws.Cells.Item(11, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(11, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(11, 19).Value = 1475
$ws.Cells.Item(11, 20).Value = 20

# Row 12
$ws.Cells.Item(12, 4).Value = 45028
$ws.Cells.Item(12, 13).Value = 200
$ws.Cells.Item(12, 14).Value = 21000
$ws.Cells.Item(12, 15).Value = 22000
$ws.Cells.Item(12, 16).Value = 21500
$ws.Cells.Item(12, 19).Value = 1075

# Row 13
$ws.Cells.Item(13, 4).Value = 44664
$ws.Cells.Item(13, 12).Value = "Segunda"
$ws.Cells.Item(13, 13).Value = 150
$ws.Cells.Item(13, 14).Value = 29000
$ws.Cells.Item(13, 15).Value = 30000
$ws.Cells.Item(13, 16).Value = 29500
$ws.Cells.Item(13, 19).Value = 1639

# Row 14
$ws.Cells.Item(14, 4).Value = 44636
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 200
$ws.Cells.Item(14, 14).Value = 29000
$ws.Cells.Item(14, 15).Value = 30000
$ws.Cells.Item(14, 16).Value = 29500
$ws.Cells.Item(14, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(14, 19).Value = 1475
$ws.Cells.Item(14, 20).Value = 20

# Row 15
$ws.Cells.Item(15, 4).Value = 45007
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 14).Value = 27000
$ws.Cells.Item(15, 15).Value = 28000
$ws.Cells.Item(15, 16).Value = 27500
$ws.Cells.Item(15, 19).Value = 1375

# Row 16
$ws.Cells.Item(16, 4).Value = 44972
$ws.Cells.Item(16, 13).Value = 140
$ws.Cells.Item(16, 14).Value = 27000
$ws.Cells.Item(16, 15).Value = 28000
$ws.Cells.Item(16, 16).Value = 27429
$ws.Cells.Item(16, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(16, 18).Value = "Región Metropolitana"
$ws.Cells.Item(16, 19).Value = 1524
$ws.Cells.Item(16, 20).Value = 18

# Row 17
$ws.Cells.Item(17, 4).Value = 45014
$ws.Cells.Item(17, 12).Value = "Segunda"
$ws.Cells.Item(17, 14).Value = 24000
$ws.Cells.Item(17, 15).Value = 25000
$ws.Cells.Item(17, 16).Value = 24500
$ws.Cells.Item(17, 19).Value = 1225

# Row 18
$ws.Cells.Item(18, 4).Value = 44979
$ws.Cells.Item(18, 13).Value = 250
$ws.Cells.Item(18, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(18, 19).Value = 1475
$ws.Cells.Item(18, 20).Value = 20
